$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.344.31'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.504.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.95'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.524'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.42%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.15'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.82%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.40'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.895.70'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.504.05'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.855'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.270.50'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.84'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0943'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.71'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +13.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.48'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '247.75'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.61'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.08'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.05'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.43'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.137'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.80'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.09'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0791'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.80%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.70'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.00'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.55%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '121.24'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.22'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0298'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.993.78'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.08'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.54%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.86%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.85'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.29%  '
